# Replace the "Name+Surname" participant-name column on the "Follow Up"
# sheet with an anonymised numeric "Code" column (1..54, matching row
# order), and update the matching entry on the "Codes" legend sheet from
# "Name+Surname" / "Name and surname" to "Code" / "Code of paticipans".

$wb = $excel.ActiveWorkbook

$followUp = $wb.Worksheets.Item("Follow Up")
$codes = $wb.Worksheets.Item("Codes")

# --- "Follow Up" sheet -----------------------------------------------
# Header
$followUp.Range("C1").Value = "Code"

# Data rows 2..55 previously held the participant's full name; replace
# with a sequential numeric code (row 2 -> 1, row 3 -> 2, ... row 55 -> 54)
for ($row = 2; $row -le 55; $row++) {
    $cell = $followUp.Cells.Item($row, 3)
    $cell.Value = $row - 1
    $cell.NumberFormat = "0"
}

# --- "Codes" legend sheet ---------------------------------------------
$codes.Range("A5").Value = "Code"
$codes.Range("B5").Value = "Code of paticipans"

# --- Selection / view state (matches the saved workbook's last-used
# selection on each sheet) ---------------------------------------------
$null = $followUp.Range("C2:C55").Select()
$null = $codes.Range("B5").Select()
$null = $followUp.Activate()
